$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data scraped by the GitHub Actions job.
# A leading apostrophe forces Excel to store the numeric-looking values as plain text
# (matching the original inlineStr cell type), and resetting the style to "Normal"
# avoids picking up an unwanted @ (Text) number-format style.

# Row 2
$ws.Range("D2").Value = "'25.928.35"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +5.48%  "
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'1.715.92"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +3.63%  "
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("D4").Value = "'1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.24%  "
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'330.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +3.37%  "
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("D6").Value = "'0.9983"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.04%  "
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("D7").Value = "'0.3687"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +1.10%  "
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("D8").Value = "'49.86"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +6.80%  "
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("D9").Value = "'0.3322"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +1.91%  "
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("D10").Value = "'1.180"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +4.49%  "
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("D11").Value = "'0.07476"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +6.13%  "
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("D12").Value = "'1.002"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.32%  "
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("D13").Value = "'6.243"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +4.45%  "
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("D14").Value = "'20.11"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +3.32%  "
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("D15").Value = "'6.909"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +4.52%  "
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("D16").Value = "'1.719.18"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +3.81%  "
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("D17").Value = "'0.00001076"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +2.89%  "
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("D18").Value = "'0.06643"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.36%  "
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("D19").Value = "'81.97"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +3.92%  "
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("D20").Value = "'0.9987"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.16%  "
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("D21").Value = "'16.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +3.58%  "
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("D22").Value = "'6.084"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +2.24%  "
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("D23").Value = "'13.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +3.57%  "
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("D24").Value = "'25.927.91"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +5.51%  "
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("D25").Value = "'2.465"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.32%  "
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("D26").Value = "'2.486"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +5.38%  "
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("D27").Value = "'150.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +1.74%  "
$ws.Range("E27").Style = "Normal"

# Row 28
$ws.Range("D28").Value = "'19.29"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +3.60%  "
$ws.Range("E28").Style = "Normal"

# Row 29
$ws.Range("D29").Value = "'1.304"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +6.91%  "
$ws.Range("E29").Style = "Normal"

# Row 30
$ws.Range("E30").Value = "'  +3.88%  "
$ws.Range("E30").Style = "Normal"

# Row 31
$ws.Range("D31").Value = "'128.89"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +3.32%  "
$ws.Range("E31").Style = "Normal"

# Row 32
$ws.Range("D32").Value = "'4.122"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +1.33%  "
$ws.Range("E32").Style = "Normal"

# Row 33
$ws.Range("D33").Value = "'5.963"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +2.92%  "
$ws.Range("E33").Style = "Normal"

# Row 34
$ws.Range("D34").Value = "'0.08543"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +0.87%  "
$ws.Range("E34").Style = "Normal"

# Row 35
$ws.Range("D35").Value = "'1.718"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +2.33%  "
$ws.Range("E35").Style = "Normal"

# Row 36
$ws.Range("D36").Value = "'12.94"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +5.30%  "
$ws.Range("E36").Style = "Normal"

# Row 37
$ws.Range("D37").Value = "'5.358"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +2.75%  "
$ws.Range("E37").Style = "Normal"

# Row 38
$ws.Range("D38").Value = "'1.284"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.93%  "
$ws.Range("E38").Style = "Normal"

# Row 39
$ws.Range("D39").Value = "'0.06223"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +3.15%  "
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("D40").Value = "'0.02289"
$ws.Range("D40").Style = "Normal"

# Row 41
$ws.Range("B41").Value = "'Algorand"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'0.2138"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +3.37%  "
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("B42").Value = "'FraxShare"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'8.552"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +5.12%  "
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("D43").Value = "'14.53"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +14.99%  "
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("D44").Value = "'0.6161"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +4.29%  "
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("D45").Value = "'0.9983"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.10%  "
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("D46").Value = "'3.834"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.44%  "
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("D47").Value = "'0.5869"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +4.43%  "
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("D48").Value = "'127.21"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +2.55%  "
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("D49").Value = "'2.012"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +3.08%  "
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("D50").Value = "'0.07256"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +4.27%  "
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("D51").Value = "'77.12"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +3.39%  "
$ws.Range("E51").Style = "Normal"

